# Update the 25 "two-digit x two-digit" practice answers in the single
# table on the page. Cells are addressed directly by (row, column) so that
# replacements which happen to collide textually with each other (e.g. a
# new value equal to some other cell's old value) can never cross-match,
# which a global Find/Replace("Replace All") pass could not guarantee.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "88×50=4400"   # was: 86×12=1032
$t.Cell(1, 2).Range.Text = "25×52=1300"   # was: 88×88=7744
$t.Cell(1, 3).Range.Text = "96×60=5760"   # was: 83×29=2407
$t.Cell(1, 4).Range.Text = "28×90=2520"   # was: 61×36=2196
$t.Cell(1, 5).Range.Text = "92×87=8004"   # was: 76×87=6612

$t.Cell(5, 1).Range.Text = "94×22=2068"   # was: 75×78=5850
$t.Cell(5, 2).Range.Text = "44×13=572"    # was: 46×25=1150
$t.Cell(5, 3).Range.Text = "73×33=2409"   # was: 71×25=1775
$t.Cell(5, 4).Range.Text = "12×86=1032"   # was: 86×21=1806
$t.Cell(5, 5).Range.Text = "58×78=4524"   # was: 26×33=858

$t.Cell(10, 1).Range.Text = "57×81=4617"  # was: 21×89=1869
$t.Cell(10, 2).Range.Text = "42×78=3276"  # was: 12×73=876
$t.Cell(10, 3).Range.Text = "84×99=8316"  # was: 57×42=2394
$t.Cell(10, 4).Range.Text = "41×94=3854"  # was: 32×60=1920
$t.Cell(10, 5).Range.Text = "32×45=1440"  # was: 15×11=165

$t.Cell(15, 1).Range.Text = "81×59=4779"  # was: 27×11=297
$t.Cell(15, 2).Range.Text = "69×50=3450"  # was: 56×14=784
$t.Cell(15, 3).Range.Text = "48×76=3648"  # was: 26×76=1976
$t.Cell(15, 4).Range.Text = "70×11=770"   # was: 37×57=2109
$t.Cell(15, 5).Range.Text = "57×70=3990"  # was: 11×27=297

$t.Cell(20, 1).Range.Text = "56×84=4704"  # was: 19×65=1235
$t.Cell(20, 2).Range.Text = "52×46=2392"  # was: 27×94=2538
$t.Cell(20, 3).Range.Text = "86×21=1806"  # was: 23×63=1449
$t.Cell(20, 4).Range.Text = "67×36=2412"  # was: 17×11=187
$t.Cell(20, 5).Range.Text = "27×31=837"   # was: 84×48=4032
